$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1500
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H33").Value = 5200.6
$ws.Range("I33").Value = 7283
$ws.Range("J33").Value = 341.66666
$ws.Range("K33").Value = 7283
$ws.Range("L33").Value = 341.66666
$ws.Range("M33").Value = -7054
$ws.Range("N33").Value = -799.66666
$ws.Range("H87").Value = 18036.264
$ws.Range("J87").Value = 18036.264
$ws.Range("L87").Value = 18036.264
$ws.Range("N87").Value = -20532.264
$ws.Range("H90").Value = 18036.264
$ws.Range("J90").Value = 18036.264
$ws.Range("L90").Value = 54108.792
$ws.Range("N90").Value = -66588.792
$ws.Range("H98").Value = 1335
$ws.Range("I98").Value = 1336
$ws.Range("K98").Value = 1336
$ws.Range("M98").Value = 162
$ws.Range("H122").Value = 1335
$ws.Range("I122").Value = 1336
$ws.Range("K122").Value = 4008
$ws.Range("M122").Value = -1558
$ws.Range("H137").Value = 2112.6667
$ws.Range("I137").Value = 1490.2142
$ws.Range("K137").Value = 4470.642599999999
$ws.Range("M137").Value = -1920.642599999999
$ws.Range("H138").Value = 1821.0488
$ws.Range("I138").Value = 1626.8948
$ws.Range("J138").Value = 1988.7273
$ws.Range("K138").Value = 4880.6844
$ws.Range("L138").Value = 5966.1819
$ws.Range("M138").Value = 259.3155999999999
$ws.Range("N138").Value = -16246.1819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3632.4524
$ws.Range("I32").Value = 2288.4805
$ws.Range("K32").Value = 2288.4805
$ws.Range("M32").Value = -2001.4805
$ws.Range("H61").Value = 6725.75
$ws.Range("I61").Value = 3728
$ws.Range("J61").Value = 9723.5
$ws.Range("K61").Value = 3728
$ws.Range("L61").Value = 9723.5
$ws.Range("M61").Value = -3516
$ws.Range("N61").Value = -10147.5
$ws.Range("H74").Value = 1516.7333
$ws.Range("I74").Value = 1617.6923
$ws.Range("J74").Value = 1439.5294
$ws.Range("K74").Value = 1617.6923
$ws.Range("L74").Value = 1439.5294
$ws.Range("M74").Value = -743.6922999999999
$ws.Range("N74").Value = -3187.5294
$ws.Range("H77").Value = 1516.7333
$ws.Range("I77").Value = 1617.6923
$ws.Range("J77").Value = 1439.5294
$ws.Range("K77").Value = 8088.461499999999
$ws.Range("L77").Value = 7197.646999999999
$ws.Range("M77").Value = -3720.461499999999
$ws.Range("N77").Value = -15933.647
$ws.Range("H132").Value = 2476.3948
$ws.Range("I132").Value = 979.4761999999999
$ws.Range("J132").Value = 4325.5293
$ws.Range("K132").Value = 2938.4286
$ws.Range("L132").Value = 12976.5879
$ws.Range("M132").Value = -408.4285999999997
$ws.Range("N132").Value = -18036.5879
$ws.Range("H136").Value = 6725.75
$ws.Range("I136").Value = 3728
$ws.Range("J136").Value = 9723.5
$ws.Range("K136").Value = 11184
$ws.Range("L136").Value = 29170.5
$ws.Range("M136").Value = -8634
$ws.Range("N136").Value = -34270.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1224.8966
$ws.Range("I99").Value = 868.4211
$ws.Range("J99").Value = 1902.2
$ws.Range("K99").Value = 868.4211
$ws.Range("L99").Value = 1902.2
$ws.Range("M99").Value = 629.5789
$ws.Range("N99").Value = -4898.2
$ws.Range("H134").Value = 2074.875
$ws.Range("I134").Value = 895.0741
$ws.Range("K134").Value = 2685.2223
$ws.Range("M134").Value = -150.2223000000004

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1722.2703
$ws.Range("I31").Value = 1477.6129
$ws.Range("J31").Value = 2986.3333
$ws.Range("K31").Value = 1477.6129
$ws.Range("L31").Value = 2986.3333
$ws.Range("M31").Value = -1182.6129
$ws.Range("N31").Value = -3576.3333
$ws.Range("H34").Value = 1722.2703
$ws.Range("I34").Value = 1477.6129
$ws.Range("J34").Value = 2986.3333
$ws.Range("K34").Value = 1477.6129
$ws.Range("L34").Value = 2986.3333
$ws.Range("M34").Value = -1275.6129
$ws.Range("N34").Value = -3390.3333
$ws.Range("H58").Value = 2281.484
$ws.Range("I58").Value = 1335.5
$ws.Range("J58").Value = 3290.5334
$ws.Range("K58").Value = 1335.5
$ws.Range("L58").Value = 3290.5334
$ws.Range("M58").Value = -1132.5
$ws.Range("N58").Value = -3696.5334
$ws.Range("H132").Value = 2806.75
$ws.Range("I132").Value = 1873.5714
$ws.Range("J132").Value = 4113.2
$ws.Range("K132").Value = 5620.7142
$ws.Range("L132").Value = 12339.6
$ws.Range("M132").Value = -3090.7142
$ws.Range("N132").Value = -17399.6
$ws.Range("H134").Value = 4137.273
$ws.Range("I134").Value = 6049.2
$ws.Range("J134").Value = 2544
$ws.Range("K134").Value = 18147.6
$ws.Range("L134").Value = 7632
$ws.Range("M134").Value = -15612.6
$ws.Range("N134").Value = -12702
$ws.Range("H136").Value = 2281.484
$ws.Range("I136").Value = 1335.5
$ws.Range("J136").Value = 3290.5334
$ws.Range("K136").Value = 4006.5
$ws.Range("L136").Value = 9871.600199999999
$ws.Range("M136").Value = -1456.5
$ws.Range("N136").Value = -14971.6002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4181.4546
$ws.Range("I132").Value = 3127.1428
$ws.Range("K132").Value = 9381.428400000001
$ws.Range("M132").Value = -6851.428400000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3196.4
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 3196.4
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 3196.4
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -3786.4
$ws.Range("H27").Value = 3196.4
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 3196.4
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 3196.4
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -3410.4
$ws.Range("H136").Value = 27784394
$ws.Range("I136").Value = 6643.143
$ws.Range("J136").Value = 66673244
$ws.Range("K136").Value = 19929.429
$ws.Range("L136").Value = 200019732
$ws.Range("M136").Value = -17379.429
$ws.Range("N136").Value = -200024832
$ws.Range("H140").Value = 61981
$ws.Range("J140").Value = 61981
$ws.Range("L140").Value = 61981
$ws.Range("N140").Value = -72341

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1675
$ws.Range("I81").Value = 2640.2
$ws.Range("J81").Value = 985.5714
$ws.Range("K81").Value = 5280.4
$ws.Range("L81").Value = 1971.1428
$ws.Range("M81").Value = -4219.4
$ws.Range("N81").Value = -4093.1428
$ws.Range("H84").Value = 1675
$ws.Range("I84").Value = 2640.2
$ws.Range("J84").Value = 985.5714
$ws.Range("K84").Value = 26402
$ws.Range("L84").Value = 9855.714
$ws.Range("M84").Value = -21098
$ws.Range("N84").Value = -20463.714
$ws.Range("H132").Value = 2184.3428
$ws.Range("I132").Value = 1943.7916
$ws.Range("J132").Value = 2709.182
$ws.Range("K132").Value = 5831.3748
$ws.Range("L132").Value = 8127.545999999999
$ws.Range("M132").Value = -3301.3748
$ws.Range("N132").Value = -13187.546
$ws.Range("H136").Value = 6284.08
$ws.Range("I136").Value = 1129.2
$ws.Range("J136").Value = 9720.666999999999
$ws.Range("K136").Value = 3387.6
$ws.Range("L136").Value = 29162.001
$ws.Range("M136").Value = -837.6000000000004
$ws.Range("N136").Value = -34262.001
